$wb = $excel.ActiveWorkbook

# --- Sponza sheet (G column performance numbers for the new v1289 build) ---
$wsSponza = $wb.Worksheets.Item("Sponza")
$wsSponza.Range("G1").Value = "v1289"
$sponzaG = @(10176,10206,10194,10139,10157,10184,10201,10175,10125,10187)
for ($i = 0; $i -lt $sponzaG.Count; $i++) {
    $row = 2 + $i
    $wsSponza.Cells.Item($row, 7).Value = $sponzaG[$i]
}
$wsSponza.Range("G15").Select()

# --- ComplexMesh sheet (G column performance numbers for the new v1289 build) ---
$wsComplex = $wb.Worksheets.Item("ComplexMesh")
$wsComplex.Range("G1").Value = "v1289"
$complexG = @(7683,7657,7648,7612,7716,7690,7631,7615,7645,7643)
for ($i = 0; $i -lt $complexG.Count; $i++) {
    $row = 2 + $i
    $wsComplex.Cells.Item($row, 7).Value = $complexG[$i]
}
$wsComplex.Range("G15").Select()

Write-Output "done"
